# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (D3) and "Correspond Handback DateTime" (G3)
# values for the second row of data (de0ea274-... file) on both the zh-cn and de-de
# language report sheets, reflecting the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-24 08:58:35"
$wsZhCn.Range("G3").Value = "2016-02-24 08:59:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-24 08:58:48"
$wsDeDe.Range("G3").Value = "2016-02-24 08:59:54"
